# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "K" (strikeouts) values replacing the previous "Strike#" derived values in column G
$kValues = @{
    2  = 0
    3  = 1
    4  = 1
    5  = 1
    6  = 0
    7  = 1
    8  = 1
    9  = 1
    10 = 0
    11 = 0
    12 = 0
    13 = 1
    14 = 0
    15 = 0
    16 = 1
    17 = 0
    18 = 0
    19 = 0
    20 = 1
    21 = 0
    22 = 2
    23 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
